# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: new case counts for several
# countries/rows, the updated "last updated" timestamp, and the
# Burkina Faso / Sri Lanka row swap (their case figures traded rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (A1) ---------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 29 de Abril de 2020 a las 20:52"

# --- Row 4: Estados Unidos ----------------------------------------------
$ws.Cells.Item(4, 2).Value = 1049431   # Casos totales
$ws.Cells.Item(4, 3).Value = 13666     # Nuevos casos
$ws.Cells.Item(4, 4).Value = 144411    # Casos activos
$ws.Cells.Item(4, 5).Value = 844380    # Recuperados
$ws.Cells.Item(4, 7).Value = 1374      # Muertes hoy
$ws.Cells.Item(4, 8).Value = 60640     # Muertes

# --- Row 15: Canada -------------------------------------------------------
$ws.Cells.Item(15, 2).Value = 51231    # Casos totales
$ws.Cells.Item(15, 3).Value = 1205     # Nuevos casos
$ws.Cells.Item(15, 5).Value = 28361    # Recuperados
$ws.Cells.Item(15, 7).Value = 125      # Muertes hoy
$ws.Cells.Item(15, 8).Value = 2984     # Muertes

# --- Row 69: Uzbekistan -----------------------------------------------
$ws.Cells.Item(69, 2).Value = 2002     # Casos totales
$ws.Cells.Item(69, 3).Value = 63       # Nuevos casos
$ws.Cells.Item(69, 5).Value = 897      # Recuperados
$ws.Cells.Item(69, 7).Value = 1        # Muertes hoy
$ws.Cells.Item(69, 8).Value = 9        # Muertes

# --- Row 96: Principado de Andorra --------------------------------------
$ws.Cells.Item(96, 4).Value = 423      # Casos activos
$ws.Cells.Item(96, 5).Value = 278      # Recuperados
$ws.Cells.Item(96, 7).Value = 1        # Muertes hoy
$ws.Cells.Item(96, 8).Value = 42       # Muertes

# --- Rows 103/104: Burkina Faso & Sri Lanka swap places -----------------
# The two countries' rows trade both their names and figures: the row
# that used to be Burkina Faso becomes Sri Lanka (with Burkina Faso's old
# numbers now updated), and vice versa.
$ws.Cells.Item(103, 1).Value = "Sri Lanka"
$ws.Cells.Item(103, 2).Value = 649
$ws.Cells.Item(103, 3).Value = 30
$ws.Cells.Item(103, 4).Value = 136
$ws.Cells.Item(103, 5).Value = 506
$ws.Cells.Item(103, 6).Value = 2
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 7

$ws.Cells.Item(104, 1).Value = "Burkina Faso"
$ws.Cells.Item(104, 2).Value = 638
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 476
$ws.Cells.Item(104, 5).Value = 120
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 42
